# Update the dSF (column F) values for the rows that were re-pulled / recalculated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Explicit row -> new F value map (row number is the worksheet row, matching the XML diff)
$rowUpdates = @{
    7  = 2
    10 = 7
    12 = -2
    22 = -5
    26 = -5
    28 = -2
    29 = 6
    34 = -5
    37 = 0
    41 = -1
    45 = -2
    47 = -2
    48 = -7
    51 = -3
    52 = -2
    53 = -2
    54 = -3
    55 = -2
    56 = -4
    61 = -5
    64 = 8
    70 = 1
    71 = -4
    74 = -2
}

foreach ($row in $rowUpdates.Keys) {
    $ws.Range("F$row").Value = $rowUpdates[$row]
}
